$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range('D2').NumberFormat = "@"
$ws.Range('D2').Value = '25.938.54'
$ws.Range('D2').Style = "Normal"
$ws.Range('E2').Value = '  -0.10%  '
$ws.Range('D3').NumberFormat = "@"
$ws.Range('D3').Value = '1.637.16'
$ws.Range('D3').Style = "Normal"
$ws.Range('E3').Value = '  -0.15%  '
$ws.Range('E4').Value = '  -0.11%  '
$ws.Range('D5').NumberFormat = "@"
$ws.Range('D5').Value = '214.62'
$ws.Range('D5').Style = "Normal"
$ws.Range('E5').Value = '  -0.12%  '
$ws.Range('D6').NumberFormat = "@"
$ws.Range('D6').Value = '0.506'
$ws.Range('D6').Style = "Normal"
$ws.Range('E6').Value = '  -0.30%  '
$ws.Range('E7').Value = '  +0.06%  '
$ws.Range('D8').NumberFormat = "@"
$ws.Range('D8').Value = '0.254'
$ws.Range('D8').Style = "Normal"
$ws.Range('E8').Value = '  -0.36%  '
$ws.Range('D9').NumberFormat = "@"
$ws.Range('D9').Value = '0.0636'
$ws.Range('D9').Style = "Normal"
$ws.Range('E9').Value = '  -0.37%  '
$ws.Range('D10').NumberFormat = "@"
$ws.Range('D10').Value = '19.50'
$ws.Range('D10').Style = "Normal"
$ws.Range('E10').Value = '  -0.92%  '
$ws.Range('D11').NumberFormat = "@"
$ws.Range('D11').Value = '0.0794'
$ws.Range('D11').Style = "Normal"
$ws.Range('E11').Value = '  -0.17%  '
$ws.Range('D12').NumberFormat = "@"
$ws.Range('D12').Value = '4.25'
$ws.Range('D12').Style = "Normal"
$ws.Range('E12').Value = '  -0.21%  '
$ws.Range('D13').NumberFormat = "@"
$ws.Range('D13').Value = '1.618.98'
$ws.Range('D13').Style = "Normal"
$ws.Range('E13').Value = '  -1.22%  '
$ws.Range('D14').NumberFormat = "@"
$ws.Range('D14').Value = '0.541'
$ws.Range('D14').Style = "Normal"
$ws.Range('E14').Value = '  -0.56%  '
$ws.Range('D15').NumberFormat = "@"
$ws.Range('D15').Value = '63.32'
$ws.Range('D15').Style = "Normal"
$ws.Range('E15').Value = '  +1.16%  '
$ws.Range('D16').NumberFormat = "@"
$ws.Range('D16').Value = '0.0₃0757'
$ws.Range('D16').Style = "Normal"
$ws.Range('E16').Value = '  -0.38%  '
$ws.Range('D17').NumberFormat = "@"
$ws.Range('D17').Value = '25.991.20'
$ws.Range('D17').Style = "Normal"
$ws.Range('E17').Value = '  +0.09%  '
$ws.Range('E18').Value = '  -0.04%  '
$ws.Range('D19').NumberFormat = "@"
$ws.Range('D19').Value = '194.19'
$ws.Range('D19').Style = "Normal"
$ws.Range('E19').Value = '  +0.03%  '
$ws.Range('D20').NumberFormat = "@"
$ws.Range('D20').Value = '4.33'
$ws.Range('D20').Style = "Normal"
$ws.Range('E20').Value = '  -0.91%  '
$ws.Range('D21').NumberFormat = "@"
$ws.Range('D21').Value = '9.87'
$ws.Range('D21').Style = "Normal"
$ws.Range('E21').Value = '  -0.72%  '
$ws.Range('D22').NumberFormat = "@"
$ws.Range('D22').Value = '6.18'
$ws.Range('D22').Style = "Normal"
$ws.Range('E22').Value = '  -1.51%  '
$ws.Range('D23').NumberFormat = "@"
$ws.Range('D23').Value = '0.132'
$ws.Range('D23').Style = "Normal"
$ws.Range('E23').Value = '  +3.90%  '
$ws.Range('D24').NumberFormat = "@"
$ws.Range('D24').Value = '143.54'
$ws.Range('D24').Style = "Normal"
$ws.Range('E24').Value = '  -0.40%  '
$ws.Range('B25').Value = 'BinanceUSD'
$ws.Range('C25').Value = 'https://coinranking.com/coin/vSo2fu9iE1s0Y+binanceusd-busd'
$ws.Range('D25').NumberFormat = "@"
$ws.Range('D25').Value = '1.00'
$ws.Range('D25').Style = "Normal"
$ws.Range('E25').Value = '  -0.17%  '
$ws.Range('B26').Value = 'Toncoin'
$ws.Range('C26').Value = 'https://coinranking.com/coin/67YlI0K1b+toncoin-ton'
$ws.Range('D26').NumberFormat = "@"
$ws.Range('D26').Value = '1.77'
$ws.Range('D26').Style = "Normal"
$ws.Range('E26').Value = '  -0.88%  '
$ws.Range('D27').NumberFormat = "@"
$ws.Range('D27').Value = '6.87'
$ws.Range('D27').Style = "Normal"
$ws.Range('E27').Value = '  +0.28%  '
$ws.Range('D28').NumberFormat = "@"
$ws.Range('D28').Value = '15.52'
$ws.Range('D28').Style = "Normal"
$ws.Range('E28').Value = '  +0.10%  '
$ws.Range('E29').Value = '  -0.16%  '
$ws.Range('D30').NumberFormat = "@"
$ws.Range('D30').Value = '0.0493'
$ws.Range('D30').Style = "Normal"
$ws.Range('E30').Value = '  -1.62%  '
$ws.Range('B31').Value = 'InternetComputer(DFINITY)'
$ws.Range('C31').Value = 'https://coinranking.com/coin/aMNLwaUbY+internetcomputerdfinity-icp'
$ws.Range('D31').NumberFormat = "@"
$ws.Range('D31').Value = '3.27'
$ws.Range('D31').Style = "Normal"
$ws.Range('E31').Value = '  -0.87%  '
$ws.Range('B32').Value = 'Filecoin'
$ws.Range('C32').Value = 'https://coinranking.com/coin/ymQub4fuB+filecoin-fil'
$ws.Range('D32').NumberFormat = "@"
$ws.Range('D32').Value = '3.24'
$ws.Range('D32').Style = "Normal"
$ws.Range('E32').Value = '  +0.33%  '
$ws.Range('D33').NumberFormat = "@"
$ws.Range('D33').Value = '1.53'
$ws.Range('D33').Style = "Normal"
$ws.Range('E33').Value = '  -0.58%  '
$ws.Range('E34').Value = '  +0.61%  '
$ws.Range('D35').NumberFormat = "@"
$ws.Range('D35').Value = '0.900'
$ws.Range('D35').Style = "Normal"
$ws.Range('E35').Value = '  -0.58%  '
$ws.Range('D36').NumberFormat = "@"
$ws.Range('D36').Value = '1.125.71'
$ws.Range('D36').Style = "Normal"
$ws.Range('E36').Value = '  -1.22%  '
$ws.Range('E37').Value = '  -1.48%  '
$ws.Range('E38').Value = '  -0.75%  '
$ws.Range('D39').NumberFormat = "@"
$ws.Range('D39').Value = '0.0156'
$ws.Range('D39').Style = "Normal"
$ws.Range('E39').Value = '  -0.49%  '
$ws.Range('D40').NumberFormat = "@"
$ws.Range('D40').Value = '98.53'
$ws.Range('D40').Style = "Normal"
$ws.Range('E40').Value = '  -0.79%  '
$ws.Range('D41').NumberFormat = "@"
$ws.Range('D41').Value = '5.41'
$ws.Range('D41').Style = "Normal"
$ws.Range('E41').Value = '  -0.15%  '
$ws.Range('D42').NumberFormat = "@"
$ws.Range('D42').Value = '0.792'
$ws.Range('D42').Style = "Normal"
$ws.Range('E42').Value = '  -1.07%  '
$ws.Range('D43').NumberFormat = "@"
$ws.Range('D43').Value = '0.0₆0115'
$ws.Range('D43').Style = "Normal"
$ws.Range('E43').Value = '  +0.42%  '
$ws.Range('D44').NumberFormat = "@"
$ws.Range('D44').Value = '56.28'
$ws.Range('D44').Style = "Normal"
$ws.Range('E44').Value = '  -0.47%  '
$ws.Range('E45').Value = '  +2.63%  '
$ws.Range('E46').Value = '  -1.63%  '
$ws.Range('D47').NumberFormat = "@"
$ws.Range('D47').Value = '7.77'
$ws.Range('D47').Style = "Normal"
$ws.Range('E47').Value = '  +1.70%  '
$ws.Range('E48').Value = '  -0.64%  '
$ws.Range('E49').Value = '  +0.14%  '
$ws.Range('D50').NumberFormat = "@"
$ws.Range('D50').Value = '0.0945'
$ws.Range('D50').Style = "Normal"
$ws.Range('E50').Value = '  -1.93%  '
$ws.Range('D51').NumberFormat = "@"
$ws.Range('D51').Value = '5.50'
$ws.Range('D51').Style = "Normal"
$ws.Range('E51').Value = '  -0.63%  '
